$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws "D2" '26.072.39'
Set-TextValue $ws "E2" '  +1.47%  '

# Row 3
Set-TextValue $ws "D3" '1.763.98'
Set-TextValue $ws "E3" '  +1.16%  '

# Row 4
Set-TextValue $ws "D4" '0.9993'
Set-TextValue $ws "E4" '  -0.13%  '

# Row 5
Set-TextValue $ws "D5" '237.74'
Set-TextValue $ws "E5" '  -0.38%  '

# Row 6
Set-TextValue $ws "E6" '  -0.13%  '

# Row 7
Set-TextValue $ws "D7" '0.5262'
Set-TextValue $ws "E7" '  +4.13%  '

# Row 8
Set-TextValue $ws "D8" '0.2743'
Set-TextValue $ws "E8" '  +3.56%  '

# Row 9
Set-TextValue $ws "D9" '0.06216'
Set-TextValue $ws "E9" '  +1.51%  '

# Row 10
Set-TextValue $ws "D10" '1.771.21'
Set-TextValue $ws "E10" '  +1.42%  '

# Row 11
$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue $ws "D11" '15.97'
Set-TextValue $ws "E11" '  +4.90%  '

# Row 12
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws "D12" '0.07041'
Set-TextValue $ws "E12" '  +1.23%  '

# Row 13
Set-TextValue $ws "D13" '0.6560'
Set-TextValue $ws "E13" '  +10.06%  '

# Row 14
Set-TextValue $ws "D14" '4.508'
Set-TextValue $ws "E14" '  +0.41%  '

# Row 15
Set-TextValue $ws "D15" '78.55'
Set-TextValue $ws "E15" '  +2.57%  '

# Row 16
Set-TextValue $ws "D16" '0.9985'
Set-TextValue $ws "E16" '  -0.16%  '

# Row 17
Set-TextValue $ws "D17" '0.9995'
Set-TextValue $ws "E17" '  -0.17%  '

# Row 18
Set-TextValue $ws "D18" '26.072.63'

# Row 19
Set-TextValue $ws "E19" '  +1.04%  '

# Row 20
Set-TextValue $ws "D20" '0.000006753'
Set-TextValue $ws "E20" '  -0.66%  '

# Row 21
Set-TextValue $ws "D21" '1.994.71'
Set-TextValue $ws "E21" '  +1.34%  '

# Row 22
Set-TextValue $ws "D22" '4.098'

# Row 23
Set-TextValue $ws "D23" '8.426'
Set-TextValue $ws "E23" '  +3.52%  '

# Row 24
Set-TextValue $ws "E24" '  +2.00%  '

# Row 25
Set-TextValue $ws "D25" '137.92'
Set-TextValue $ws "E25" '  +0.39%  '

# Row 26
Set-TextValue $ws "D26" '1.486'
Set-TextValue $ws "E26" '  -1.93%  '

# Row 27
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws "D27" '15.24'
Set-TextValue $ws "E27" '  +1.90%  '

# Row 28
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws "D28" '1.841'
Set-TextValue $ws "E28" '  +1.08%  '

# Row 29
Set-TextValue $ws "D29" '103.17'
Set-TextValue $ws "E29" '  -0.03%  '

# Row 30
Set-TextValue $ws "D30" '0.08431'
Set-TextValue $ws "E30" '  +4.08%  '

# Row 31
Set-TextValue $ws "D31" '3.713'
Set-TextValue $ws "E31" '  -1.13%  '

# Row 32
Set-TextValue $ws "D32" '3.451'
Set-TextValue $ws "E32" '  -0.24%  '

# Row 33
Set-TextValue $ws "D33" '0.04435'
Set-TextValue $ws "E33" '  -1.23%  '

# Row 34
Set-TextValue $ws "D34" '2.650'
Set-TextValue $ws "E34" '  -0.14%  '

# Row 35
Set-TextValue $ws "D35" '1.003'
Set-TextValue $ws "E35" '  +2.11%  '

# Row 36
Set-TextValue $ws "D36" '0.6107'
Set-TextValue $ws "E36" '  +0.37%  '

# Row 37
Set-TextValue $ws "D37" '2.753'
Set-TextValue $ws "E37" '  +3.86%  '

# Row 38
Set-TextValue $ws "D38" '0.01587'
Set-TextValue $ws "E38" '  +2.36%  '

# Row 39
Set-TextValue $ws "D39" '1.970'
Set-TextValue $ws "E39" '  +2.61%  '

# Row 40
Set-TextValue $ws "E40" '  +0.04%  '

# Row 41
Set-TextValue $ws "D41" '103.19'
Set-TextValue $ws "E41" '  -0.32%  '

# Row 42
Set-TextValue $ws "D42" '0.3921'
Set-TextValue $ws "E42" '  +3.27%  '

# Row 43
Set-TextValue $ws "D43" '0.7543'
Set-TextValue $ws "E43" '  +3.69%  '

# Row 44
Set-TextValue $ws "D44" '4.983'
Set-TextValue $ws "E44" '  -2.63%  '

# Row 45
Set-TextValue $ws "D45" '0.05498'
Set-TextValue $ws "E45" '  +3.14%  '

# Row 46
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws "D46" '6.284'
Set-TextValue $ws "E46" '  +6.68%  '

# Row 47
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws "D47" '0.1126'
Set-TextValue $ws "E47" '  +1.55%  '

# Row 48
Set-TextValue $ws "D48" '30.24'
Set-TextValue $ws "E48" '  +0.58%  '

# Row 49
Set-TextValue $ws "D49" '52.94'
Set-TextValue $ws "E49" '  +1.09%  '

# Row 50
$ws.Range("B50").Value = 'Decentraland'
$ws.Range("C50").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue $ws "D50" '0.3466'
Set-TextValue $ws "E50" '  +0.89%  '

# Row 51
$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
Set-TextValue $ws "D51" '1.002'
Set-TextValue $ws "E51" '  +0.37%  '
